# Refresh the cached "datetimeFigureOut" footer-date text from 10/11/2023 to
# 13/11/2023 everywhere it appears: on the Slide Master and on every one of
# its slide layouts (CustomLayouts). ppPlaceholderDate == 16 is used to find
# the date placeholder shape on each container regardless of its shape index.

$p = $ppt.ActivePresentation
$oldDate = "10/11/2023"
$newDate = "13/11/2023"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($container) {
    $shapes = $container.Shapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        $phType = $null
        try { $phType = $shape.PlaceholderFormat.Type } catch { $phType = $null }
        if ($phType -eq $ppPlaceholderDate -and $shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide Master
Update-DatePlaceholder $p.SlideMaster

# Every slide layout that hangs off the master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li)
}
